# COREESG_holdings.xlsx edit — daily model-holdings refresh.
# The worksheet carries (legacy-hashed) protection, so we must unprotect,
# apply the cell/text updates, then re-protect before handing the workbook
# back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Bump the "as of" date in the confidential disclaimer banner (A10).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) for each fund row.
$ws.Range("D2").Value = 0.253384121976319
$ws.Range("E2").Value = -0.006951180084060704

$ws.Range("D3").Value = 0.4909409652944576
$ws.Range("E3").Value = -0.007203842049092746

$ws.Range("D4").Value = 0.1004466846211413
$ws.Range("E4").Value = 0.0009433962264149276

$ws.Range("D5").Value = 0.09892810579484773
$ws.Range("E5").Value = -0.005218070709043676

$ws.Range("D6").Value = 0.0563001223132344
$ws.Range("E6").Value = -0.003176764238711161

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = -0.005898284874830795

$ws.Protect()
